# Fix column header name in the BOM table: "JCSC Part #" -> "LCSC Part #"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the mis-typed column header (this also renames the matching
# table column / query-table field, since the header cell drives the
# ListObject column name).
$ws.Range("E1").Value = "LCSC Part #"

# Leave the selection on the first data cell below the renamed header,
# matching where the user would land after editing the header and
# pressing Enter.
[void]$ws.Range("E2").Select()
